# Insert a new data row at row 213 (pushing the existing rows 213-333 down to 214-334)
# and populate it with the new observation's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("213:213").Insert()

$ws.Cells.Item(213, 1).Value  = 10
$ws.Cells.Item(213, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(213, 3).Value  = "La Araucanía"
$ws.Cells.Item(213, 4).Value  = 45176
$ws.Cells.Item(213, 5).Value  = 9
$ws.Cells.Item(213, 6).Value  = 100112005
$ws.Cells.Item(213, 7).Value  = "Puerro"
$ws.Cells.Item(213, 8).Value  = "Azul de Maquehue"
$ws.Cells.Item(213, 9).Value  = "Primera"
$ws.Cells.Item(213, 10).Value = 50
$ws.Cells.Item(213, 11).Value = 9000
$ws.Cells.Item(213, 12).Value = 9000
$ws.Cells.Item(213, 13).Value = 9000
$ws.Cells.Item(213, 14).Value = "`$/docena de paquetes"
$ws.Cells.Item(213, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(213, 16).Value = 750
$ws.Cells.Item(213, 17).Value = 12
$ws.Cells.Item(213, 18).Value = "Hortaliza"
